# "files updated and bug fixed"
# Corrected figures for rows 4 ("registered families") and 5 ("families
# receiving subsistence allowance") for years 2015-2018 (columns E-H) in
# the ზუგდიდი municipality sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - რეგისტრირებული ოჯახი (registered family)
$ws.Range("E4").Value = 15994
$ws.Range("F4").Value = 13405
$ws.Range("G4").Value = 9689
$ws.Range("H4").Value = 9711

# Row 5 - საარსებო შემწეობის მიმღები ოჯახი (family receiving subsistence allowance)
$ws.Range("E5").Value = 3676
$ws.Range("F5").Value = 4562
$ws.Range("G5").Value = 4294
$ws.Range("H5").Value = 4707
